$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A79").NumberFormat = "@"
$ws.Range("A79").Value = "2025-10-15"
$ws.Range("B79").Value = "15:21:44"
$ws.Range("C79").Value = "1.00 EUR = 1,675.9082"
